# Append a new data row (row 3) to the worksheet, mirroring the schema
# used by the existing header (row 1) / data (row 2) rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric cells ---
$ws.Range("A3").Value = 112334575
$ws.Range("B3").Value = 89539
$ws.Range("E3").Value = 1202
$ws.Range("Q3").Value = 621283
$ws.Range("R3").Value = 6614833
$ws.Range("S3").Value = 25

# --- Text cells ---
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "NT"
$ws.Range("F3").Value = "Ullticka"
$ws.Range("G3").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H3").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("P3").Value = "skjutbanan, Upl"
$ws.Range("T3").Value = "Uppsala"
$ws.Range("U3").Value = "Enköping"
$ws.Range("V3").Value = "Uppland"
$ws.Range("W3").Value = "Litslena"
$ws.Range("AW3").Value = "Elin Sjögren Englund"
$ws.Range("AX3").Value = "Elin Sjögren Englund"

# --- Date-looking cells that must stay plain text (leading apostrophe
#     forces text interpretation instead of an Excel date serial) ---
$ws.Range("Y3").Value = "'2023-09-24"
$ws.Range("AA3").Value = "'2023-09-24"

# --- Boolean cells ---
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false

# --- Present-but-empty placeholder cells (match the blank, typed cells
#     that already exist for these columns on other rows) ---
$ws.Range("I3").Value = "'"
$ws.Range("J3").Value = "'"
$ws.Range("K3").Value = "'"
$ws.Range("N3").Value = "'"
$ws.Range("AF3").Value = "'"
$ws.Range("AT3").Value = "'"
$ws.Range("AY3").Value = "'"
